# Refresh the crypto price/volume table with the latest scraped values.
# (GitHub Actions scheduled update -- see commit message.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '70.817.57'
$ws.Cells.Item(2, 5).Value = '  -0.78%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.572.08'
$ws.Cells.Item(3, 5).Value = '  -2.05%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''583.85'
$ws.Cells.Item(5, 5).Value = '  +0.25%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''184.10'
$ws.Cells.Item(6, 5).Value = '  -2.90%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.560.86'
$ws.Cells.Item(7, 5).Value = '  -2.18%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.621'
$ws.Cells.Item(8, 5).Value = '  -1.23%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.01%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.214'
$ws.Cells.Item(10, 5).Value = '  +13.05%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.647'
$ws.Cells.Item(11, 5).Value = '  -2.44%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''54.11'
$ws.Cells.Item(12, 5).Value = '  -1.89%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.23%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''9.49'
$ws.Cells.Item(14, 5).Value = '  -2.82%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.131.50'
$ws.Cells.Item(15, 5).Value = '  -2.21%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '70.691.85'
$ws.Cells.Item(16, 5).Value = '  -0.69%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''19.31'
$ws.Cells.Item(17, 5).Value = '  -2.80%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.569.55'
$ws.Cells.Item(18, 5).Value = '  -1.82%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(19, 4).Value = '''12.39'
$ws.Cells.Item(19, 5).Value = '  -2.33%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).Value = '''568.22'
$ws.Cells.Item(20, 5).Value = '  +12.40%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.60%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -5.75%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''17.68'
$ws.Cells.Item(23, 5).Value = '  -10.60%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''4.59'
$ws.Cells.Item(24, 5).Value = '  +3.39%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''4.99'
$ws.Cells.Item(25, 5).Value = '  -0.46%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''94.79'

# Row 27
$ws.Cells.Item(27, 4).Value = '''11.20'
$ws.Cells.Item(27, 5).Value = '  -3.85%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''2.92'
$ws.Cells.Item(28, 5).Value = '  -4.03%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''9.09'
$ws.Cells.Item(29, 5).Value = '  -4.68%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''32.21'
$ws.Cells.Item(30, 5).Value = '  -0.24%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''7.29'
$ws.Cells.Item(31, 5).Value = '  -7.05%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''12.27'
$ws.Cells.Item(32, 5).Value = '  -3.94%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''64.21'
$ws.Cells.Item(33, 5).Value = '  -3.48%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -3.19%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''3.32'
$ws.Cells.Item(35, 5).Value = '  +2.25%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''547.59'
$ws.Cells.Item(36, 5).Value = '  -5.64%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.412'
$ws.Cells.Item(37, 5).Value = '  -0.84%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.0₃0803'
$ws.Cells.Item(38, 5).Value = '  -1.16%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.32%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''37.41'
$ws.Cells.Item(40, 5).Value = '  -5.20%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '3.461.64'
$ws.Cells.Item(41, 5).Value = '  +6.44%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.136'
$ws.Cells.Item(42, 5).Value = '  -1.79%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''3.37'
$ws.Cells.Item(43, 5).Value = '  -4.17%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''3.08'
$ws.Cells.Item(44, 5).Value = '  -8.09%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''3.54'
$ws.Cells.Item(45, 5).Value = '  -0.50%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'ThetaToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(46, 4).Value = '''2.95'
$ws.Cells.Item(46, 5).Value = '  -4.88%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value = '''0.0441'
$ws.Cells.Item(47, 5).Value = '  -3.79%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''9.31'
$ws.Cells.Item(48, 5).Value = '  -3.87%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.18%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.998'
$ws.Cells.Item(50, 5).Value = '  +0.06%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''1.43'
$ws.Cells.Item(51, 5).Value = '  -5.99%  '
